$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Teczki akt osobowych"
$ws.Range("D10").Value = "Wordy\Error Teczki akt osobowych.docx"

$ws.Hyperlinks.Add($ws.Range("D10"), "Wordy\Error Teczki akt osobowych.docx") | Out-Null

$ws.Range("D10").Style = "Hiperłącze"

$ws.Range("D10").Select()
